$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time range text for row 26 (B26) and hours value (G26)
$ws.Range("B26").Value = "12.45-14.45"
$ws.Range("G26").Value = 2

# Update the selected/active cell to H26 (matches final selection in the file)
$ws.Range("H26").Select()
